# Daily Update 키워드 10개
# Applies the day's price / product refresh for the 마우스 (mouse) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price refreshes -------------------------------------------------
$ws.Range("E8").Value  = "37830"
$ws.Range("E42").Value = "24860"
$ws.Range("E55").Value = "23480"
$ws.Range("E60").Value = "52330"
$ws.Range("E67").Value = "12570"

# --- Rows 28 / 29: listings swapped ------------------------------------------
$ws.Range("B28").Value = "로지텍G G302 Daedalus Prime"
$ws.Range("C28").Value = "https://search.shopping.naver.com/gate.nhn?id=8134705034"
$ws.Range("D28").Value = "https://shopping-phinf.pstatic.net/main_8134705/8134705034.20211210121403.jpg"
$ws.Range("I28").Value = "로지텍G"
$ws.Range("J28").Value = "로지텍"

$ws.Range("B29").Value = "세컨드찬스 긱스타 GM900 3325 LED 게이밍 마우스"
$ws.Range("C29").Value = "https://search.shopping.naver.com/gate.nhn?id=19187823193"
$ws.Range("D29").Value = "https://shopping-phinf.pstatic.net/main_1918782/19187823193.20200612120756.jpg"
$ws.Range("I29").Value = "긱스타"
$ws.Range("J29").Value = "세컨드찬스"

# --- Rows 51 / 52 / 53: three-way rotation ------------------------------------
$ws.Range("B51").Value = "로지텍 G PRO 게이밍 마우스"
$ws.Range("C51").Value = "https://search.shopping.naver.com/gate.nhn?id=11309020076"
$ws.Range("D51").Value = "https://shopping-phinf.pstatic.net/main_1130902/11309020076.20170406170430.jpg"
$ws.Range("E51").Value = "69350"
$ws.Range("I51").Value = "로지텍"
$ws.Range("J51").Value = "로지텍"

$ws.Range("B52").Value = "로지텍 M90"
$ws.Range("C52").Value = "https://search.shopping.naver.com/gate.nhn?id=6777128876"
$ws.Range("D52").Value = "https://shopping-phinf.pstatic.net/main_6777128/6777128876.20210618110511.jpg"
$ws.Range("E52").Value = "4910"
$ws.Range("I52").Value = "로지텍"
$ws.Range("J52").Value = "로지텍"

$ws.Range("B53").Value = "삼성전자 삼성 SPA-MMG1PUB 게이밍마우스"
$ws.Range("C53").Value = "https://search.shopping.naver.com/gate.nhn?id=22366713052"
$ws.Range("D53").Value = "https://shopping-phinf.pstatic.net/main_2236671/22366713052.20200527173225.jpg"
$ws.Range("E53").Value = "29900"
$ws.Range("I53").Value = "삼성"
$ws.Range("J53").Value = "삼성전자"

# --- Rows 62 / 63: listings swapped -------------------------------------------
$ws.Range("B62").Value = "플레오맥스 MO-ER700 인체공학 버티칼 마우스"
$ws.Range("C62").Value = "https://search.shopping.naver.com/gate.nhn?id=26232573526"
$ws.Range("D62").Value = "https://shopping-phinf.pstatic.net/main_2623257/26232573526.20210304120706.jpg"
$ws.Range("E62").Value = "19000"
$ws.Range("I62").Value = "플레오맥스"
$ws.Range("J62").Value = "플레오맥스"

$ws.Range("B63").Value = "스카이디지탈 NMOUSE 4K REMASTER"
$ws.Range("C63").Value = "https://search.shopping.naver.com/gate.nhn?id=23891813522"
$ws.Range("D63").Value = "https://shopping-phinf.pstatic.net/main_2389181/23891813522.20200825100651.jpg"
$ws.Range("E63").Value = "19800"
$ws.Range("I63").Value = "스카이디지탈"
$ws.Range("J63").Value = "스카이디지탈"

# --- Rows 94 / 95: listings swapped -------------------------------------------
$ws.Range("B94").Value = "ROCCAT KONE PURE ULTRA i"
$ws.Range("C94").Value = "https://search.shopping.naver.com/gate.nhn?id=22100501353"
$ws.Range("D94").Value = "https://shopping-phinf.pstatic.net/main_2210050/22100501353.20210623094204.jpg"
$ws.Range("E94").Value = "100000"
$ws.Range("I94").Value = "ROCCAT"
$ws.Range("J94").Value = "ROCCAT"

$ws.Range("B95").Value = "지클릭커 GM-M250 LED 무소음 게이밍 마우스"
$ws.Range("C95").Value = "https://search.shopping.naver.com/gate.nhn?id=18985762773"
$ws.Range("D95").Value = "https://shopping-phinf.pstatic.net/main_1898576/18985762773.20200101140950.jpg"
$ws.Range("E95").Value = "4190"
$ws.Range("I95").Value = "지클릭커"
$ws.Range("J95").Value = "지클릭커"
